# Refresh cryptos list values (scheduled GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.58"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "1.551.63"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'206.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "'0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "1.772.69"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "1.554.06"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'3.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "26.896.53"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "'61.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "'216.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").Value = "'9.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'153.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'14.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "1.425.14"
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  +4.43%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "'0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "'5.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.986"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("D45").Value = "'63.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "1.687.54"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").Value = "'86.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("D50").Value = "0.0₇0977"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.77%  "
